# #327 Ajout des profils d'acces a58d18c1e8091c98efec92c8c093b361a253eee5
#
# 1) Update the "Date" metadata value on the Metadata sheet.
# 2) On the Elements sheet, the two "Mapping" columns (AK = RIM Mapping,
#    AL = Spécification métier) are swapped: column AK now holds what used
#    to be in AL (and vice-versa) for both the header and the data rows,
#    and the column widths follow the (now wider/narrower) content.

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B8").Value = "2024-03-19T13:17:15+00:00"

$wsElements = $wb.Worksheets.Item("Elements")

# Only rows 1, 3, 5 and 6 actually have differing AK/AL content (rows 2 and 4
# are blank in both columns), so only those need to be touched.
$rowsToSwap = @(1, 3, 5, 6)
foreach ($r in $rowsToSwap) {
    $akCell = $wsElements.Cells.Item($r, 37)
    $alCell = $wsElements.Cells.Item($r, 38)

    $akValue = $akCell.Value()
    $alValue = $alCell.Value()

    $akCell.Value = $alValue
    $alCell.Value = $akValue
}

# Column widths follow the content that now lives in each column: AK holds
# the long "Spécification métier..." text, AL holds the short "RIM Mapping".
# (values chosen as the closest achievable match to the bestFit target widths
# of 75.78515625 / 24.98046875 given the host's column-width rounding model)
$wsElements.Columns.Item(37).ColumnWidth = 75
$wsElements.Columns.Item(38).ColumnWidth = 24.166666666666668
